$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark that currently sits right after
#    "... as an output" (Word drops this automatically on save; we re-add
#    it at its new location once the surrounding content has moved).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Locate the "Android Chess Game" project entry. It is the paragraph
#    right after a blank separator paragraph, and is followed by its bullet
#    description, its "Technologies:" line, and a trailing blank paragraph.
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Android Chess Game*") {
        $titlePara = $p
        break
    }
}

$precedingBlank = $titlePara.Previous()
$blockEnd = $titlePara.Next().Next().Next()

# 3. Delete the whole "Android Chess Game" block (title, bullet,
#    technologies line, trailing blank paragraph) in one shot, leaving the
#    blank separator paragraph before it intact.
$delRange = $d.Range($titlePara.Range.Start, $blockEnd.Range.End)
$delRange.Delete()

# 4. Re-add the _GoBack bookmark, collapsed (zero-length) inside the blank
#    paragraph that now sits directly before "Udacity Android Development
#    Course".
$bmRange = $precedingBlank.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)
